$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2..63 (row 1 is header).
# We need to insert two brand-new rows right before the current row 40
# (current row 40 data will become row 42), shifting all rows 40..63 down
# to 42..65, and fill the two new rows (40 and 41) with new values.

# Insert two blank rows at position 40 (pushing existing row 40 downward twice)
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(40).Insert()

# Populate the two new rows by duplicating the row that is now at 42
# (this carries over all the common columns: A, B, C, E, F, G, H, I, J, K, R
# and also copies formatting, e.g. the date style on column D)
$ws.Range("A42:T42").Copy()
$ws.Range("A40:T40").PasteSpecial()
$ws.Range("A42:T42").Copy()
$ws.Range("A41:T41").PasteSpecial()
$excel.CutCopyMode = $false

# Now set the specific values that differ for the two new rows (40 and 41)

# Row 40 - "Primera" quality, new pricing
$ws.Cells.Item(40, 4).Value = 45072          # D - Fecha
$ws.Cells.Item(40, 12).Value = 'Primera'     # L - Calidad
$ws.Cells.Item(40, 13).Value = 250           # M - Volumen
$ws.Cells.Item(40, 14).Value = 5000          # N - Precio minimo
$ws.Cells.Item(40, 15).Value = 6000          # O - Precio maximo
$ws.Cells.Item(40, 16).Value = 5600          # P - Precio promedio ponderado
$ws.Cells.Item(40, 17).Value = '$/caja 10 kilos'  # Q - Unidad de comercializacion
$ws.Cells.Item(40, 19).Value = 560           # S - Precio $/Kg
$ws.Cells.Item(40, 20).Value = 10            # T - Kg / unidad

# Row 41 - "Segunda" quality, new pricing
$ws.Cells.Item(41, 4).Value = 45072          # D - Fecha
$ws.Cells.Item(41, 12).Value = 'Segunda'     # L - Calidad
$ws.Cells.Item(41, 13).Value = 150           # M - Volumen
$ws.Cells.Item(41, 14).Value = 4000          # N - Precio minimo
$ws.Cells.Item(41, 15).Value = 5000          # O - Precio maximo
$ws.Cells.Item(41, 16).Value = 4467          # P - Precio promedio ponderado
$ws.Cells.Item(41, 17).Value = '$/caja 10 kilos'  # Q - Unidad de comercializacion
$ws.Cells.Item(41, 19).Value = 447           # S - Precio $/Kg
$ws.Cells.Item(41, 20).Value = 10            # T - Kg / unidad
